$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(9, 3).Value = '[name="L.G.D. Officer"]   Do you copy? All targets in Sector 7 are confirmed to be ''Terracotta''!
'
$ws.Cells.Item(46, 3).Value = '[name="Ch''en"]   ''Terracotta.'' A fitting name.
'
$ws.Cells.Item(47, 3).Value = '[name="Wei Fumizuki"]   Calling them ''Reunion'' would cause unnecessary panic, Officer Ch''en. They are not flesh and blood people.
'
$ws.Cells.Item(50, 3).Value = '[name="Wei Fumizuki"]   We''ve been unable to discern the true nature of the ''Terracottas,'' but their strength has continued to grow over the past few hours.
'
$ws.Cells.Item(59, 3).Value = '[name="Wei Fumizuki"]   Our conventional methods have failed. We have no indication of how the enemy ''appeared.''
'
$ws.Cells.Item(119, 3).Value = '[name="Wei Fumizuki"]   The leader of our partner organization, Rhodes Island Pharmaceuticals. I believe we can put it that way for now. ''Doctor'' is the appropriate form of address.
'
$ws.Cells.Item(121, 3).Value = '[name="L.G.D. Officer"]   All information on the ''Terracottas'' has been disseminated to all members of the L.G.D. through the terminal.
'
$ws.Cells.Item(126, 3).Value = '[name="L.G.D. Officer"]   We can''t confirm that the ''Terracottas'' are directly connected to the heat source. And it is likely that the defensive line will collapse if we reallocate manpower.
'
$ws.Cells.Item(141, 3).Value = '[name="Wei Fumizuki"]   Nian sightings have become less and less common since the development of nomadic city technology. We have not had the opportunity to study the ''Nian''.
'
$ws.Cells.Item(142, 3).Value = '[name="Wei Fumizuki"]   The L.G.D. will do whatever it takes to uncover the secret of the ''Nian,'' even if it''s just the tip of the iceberg.
'
$ws.Cells.Item(210, 3).Value = '[name="Back-Alley Doctor"]   I think you''ve got it a bit twisted, Waai Fu. I''m afraid these aren''t quite ''people,'' you know?
'
$ws.Cells.Item(258, 3).Value = '[name="Back-Alley Doctor"]   By the way, Hung, aren''t you curious about the ''Nian''?
'
$ws.Cells.Item(284, 3).Value = '[name="???"]   Isn''t a sweet little show like that supposed to make you think, ''Oh, maybe she''s not so bad,'' and lower the dagger?
'
$ws.Cells.Item(292, 3).Value = '[name="Lava"]   Besides, there''s no way you could understand us, ''Nian''.
'
$ws.Cells.Item(294, 3).Value = '[name="Nian"]   That''s what they used to call me, but I don''t think ''Nian'' is cool enough. Can we change it? I want a badass codename like ''Lava.''
'
$ws.Cells.Item(330, 3).Value = '[name="Nian"]   Is this like a ''surrender now and I''ll let you live'' thing?
'
$ws.Cells.Item(367, 3).Value = '[name="Nian"]   But I do have one teensy little question. How did you come up with all this stuff while I, or I guess ''we,'' were snoozing away?
'
$ws.Cells.Item(394, 3).Value = '[name="Lava"]   ''The heart of the crucible is... nothing but straw!''
'
$ws.Cells.Item(404, 3).Value = '[name="Lava"]   Old folk sayings like ''the din of fireworks rings in the new year'' and such.
'
$ws.Cells.Item(419, 3).Value = '[name="Lava"]   So this is the power of the ''Nian,'' huh? I''m not impressed.
'
$ws.Cells.Item(431, 3).Value = '[name="Nian"]   ''Heavenly forge, take the bounty of the earth and smelt an edge to slice the sky!''
'

$ws.Cells.Item(94, 4).Value = '[name="오니 누님"]   아니, 그거 말고. 너희 근위국 놈들은 평소에 대체 어떤 놈들이랑 나다니는 거냐? 우르수스 사람? 아니면 강제 사람?
'
